# Parametrização de pesos de carga, spread e pico máximo de turmas permitido
#
# Updates the consolidated instructor/project allocation table (B2:J19)
# with the recalculated distribution values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..J (DD1, DD2_Onda1, DD2_Onda2, DD2_Onda3, IT1,
# IT2_Onda1, IT2_Onda2, IT2_Onda3, Total), one row per instructor (rows 2..19).
$data = @(
    @(0, 1, 2, 2, 2, 3, 1, 3, 14),  # PROG_1
    @(1, 1, 2, 1, 1, 1, 3, 4, 14),  # PROG_10
    @(0, 1, 3, 1, 1, 0, 3, 4, 13),  # PROG_11
    @(1, 5, 4, 1, 1, 1, 1, 0, 14),  # PROG_12
    @(2, 2, 2, 2, 0, 4, 0, 2, 14),  # PROG_13
    @(0, 1, 0, 3, 2, 3, 3, 2, 14),  # PROG_2
    @(0, 4, 2, 2, 2, 1, 2, 1, 14),  # PROG_3
    @(0, 1, 2, 3, 2, 1, 3, 2, 14),  # PROG_4
    @(0, 1, 2, 3, 2, 1, 3, 2, 14),  # PROG_5
    @(0, 1, 2, 2, 1, 3, 1, 3, 13),  # PROG_6
    @(1, 2, 1, 1, 2, 2, 3, 2, 14),  # PROG_7
    @(2, 4, 0, 0, 2, 3, 2, 1, 14),  # PROG_8
    @(1, 0, 2, 5, 0, 3, 1, 1, 13),  # PROG_9
    @(0, 1, 1, 1, 0, 1, 4, 4, 12),  # ROB_1
    @(0, 3, 1, 1, 0, 1, 2, 3, 11),  # ROB_2
    @(0, 2, 2, 3, 0, 2, 1, 1, 11),  # ROB_3
    @(0, 2, 3, 3, 0, 2, 0, 1, 11),  # ROB_4
    @(0, 2, 3, 3, 0, 2, 1, 1, 12)   # ROB_5
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 2 + $j   # column B = 2 .. J = 10
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
